{"js": "// Insert a new, empty paragraph right after the paragraph that contains\n// \"Databasen skal designes s\u00e5 den opfylder 3. normalform.\" \u2014 the new\n// paragraph inherits its formatting (pPr/rPr) from that paragraph, same\n// as Word does when splitting/adding a paragraph at that location.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"Databasen skal designes s\u00e5 den opfylder 3. normalform.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target paragraph text not found.\");\n}\n\nconst targetParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// Insert an empty paragraph immediately after the target paragraph.\ntargetParagraph.insertParagraph(\"\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new, empty paragraph right after the paragraph that contains\n# \"Databasen skal designes s\u00e5 den opfylder 3. normalform.\" \u2014 the new\n# paragraph inherits its formatting (pPr/rPr) from that paragraph, same\n# as Word does when inserting a paragraph break at that location.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"Databasen skal designes s\u00e5 den opfylder 3. normalform.\")\n\nif (-not $found) {\n    throw \"Target paragraph text not found.\"\n}\n\n$range.InsertParagraphAfter()\n"}
